$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 36164
$ws.Range("G2").Value = 2047
$ws.Range("H2").Value = 53
$ws.Range("I2").Value = 55

# Row 3
$ws.Range("B3").Value = "Dartin Dan"
$ws.Range("C3").Value = 31
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 19917
$ws.Range("G3").Value = 1150
$ws.Range("H3").Value = 51.96
$ws.Range("I3").Value = 34
$ws.Range("J3").Value = 0

# Row 4
$ws.Range("B4").Value = "Rocky Van Den Eeckhoudt"
$ws.Range("C4").Value = 27
$ws.Range("D4").Value = 3
$ws.Range("F4").Value = 16278
$ws.Range("G4").Value = 888
$ws.Range("H4").Value = 54.99
$ws.Range("I4").Value = 30
$ws.Range("J4").Value = 1

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Gijs Tromp"
$ws.Range("C5").Value = 27
$ws.Range("F5").Value = 16912
$ws.Range("G5").Value = 920
$ws.Range("H5").Value = 55.15
$ws.Range("I5").Value = 29

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Alessandro Delia"
$ws.Range("C6").Value = 26
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 19418
$ws.Range("G6").Value = 1138
$ws.Range("H6").Value = 51.19

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Nick Fitzpatrick"
$ws.Range("C7").Value = 24
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 11558
$ws.Range("G7").Value = 544
$ws.Range("H7").Value = 63.74
$ws.Range("I7").Value = 26
$ws.Range("J7").Value = 1

# Row 8
$ws.Range("B8").Value = "Robin Willis"
$ws.Range("C8").Value = 22
$ws.Range("F8").Value = 17535
$ws.Range("G8").Value = 1057
$ws.Range("H8").Value = 49.77
$ws.Range("I8").Value = 22

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Max Walter"
$ws.Range("C9").Value = 18
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 9206
$ws.Range("G9").Value = 663
$ws.Range("H9").Value = 41.66

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Louis Tweddle"
$ws.Range("C10").Value = 17
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 11864
$ws.Range("G10").Value = 581
$ws.Range("H10").Value = 61.26
$ws.Range("I10").Value = 18

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Noah B"
$ws.Range("C17").Value = 7
$ws.Range("F17").Value = 6577
$ws.Range("G17").Value = 380
$ws.Range("H17").Value = 51.92
$ws.Range("I17").Value = 7

# Row 18
$ws.Range("B18").Value = "Tristan Snoep"
$ws.Range("C18").Value = 6
$ws.Range("F18").Value = 7796
$ws.Range("G18").Value = 573
$ws.Range("H18").Value = 40.82
$ws.Range("I18").Value = 6

# Row 19
$ws.Range("A19").Value = 18

